$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = -10
$ws.Range("B5").Value = 35
$ws.Range("B6").Value = 35
$ws.Range("B7").Value = 70
$ws.Range("B8").Value = "EU23"
